# Fill in the missing "week 23" entries that were left blank in the
# logboek/urenoverzicht workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Logboek": add the work-log description for week 23 (row 23) ---
$ws2 = $wb.Worksheets.Item("Logboek")
$ws2.Activate()
$ws2.Range("B23").Value = "gewerkt aan de documentatie te verbeteren en in te voegen"
[void]$ws2.Range("B24").Select()

# --- Sheet "Uren periode 1 en 2": fill in the hours worked for week 23 (row 23) ---
$ws1 = $wb.Worksheets.Item("Uren periode 1 en 2")
$ws1.Activate()

$ws1.Range("B23").Value = 3
$ws1.Range("C23").Value = 24
$ws1.Range("E23").Value = 24
$ws1.Range("F23").Value = "SH akkoord"

# Leave the selection where the edit happened, on the sheet that stays active.
[void]$ws1.Range("P13").Select()

$excel.Calculate()

Write-Output "edit applied"
